$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ibuprofen")

# Update the RxCUI code for Ibuprofen from the text "C0020740" to the
# correct numeric RxCUI value 5640.
$ws.Range("D2").Value = 5640
[void]$ws.Range("D2").Select()
